$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in newly-computed experiment values
$ws.Range("E10").Value = 0.048639206997117899
$ws.Range("F11").Value = 0.97811823366959705

# Update the active selection to match the author's final cursor position
$ws.Range("E10").Select()
